$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.746.62"
$ws.Range("E2").Value = "  -0.45%  "
$ws.Range("D3").Value = "1.632.11"
$ws.Range("E3").Value = "  -0.59%  "
$ws.Range("D5").Value = "'214.75"
$ws.Range("E5").Value = "  -0.53%  "
$ws.Range("D6").Value = "'0.502"
$ws.Range("E6").Value = "  -0.90%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("E8").Value = "  -0.64%  "
$ws.Range("E9").Value = "  -0.99%  "
$ws.Range("D10").Value = "'19.55"
$ws.Range("E10").Value = "  -4.91%  "
$ws.Range("D11").Value = "'0.0784"
$ws.Range("E11").Value = "  +0.14%  "
$ws.Range("D12").Value = "1.640.90"
$ws.Range("E12").Value = "  -0.05%  "
$ws.Range("D13").Value = "'4.24"
$ws.Range("E13").Value = "  -1.09%  "
$ws.Range("D14").Value = "1.857.39"
$ws.Range("E14").Value = "  -0.63%  "
$ws.Range("D15").Value = "'0.553"
$ws.Range("E15").Value = "  -1.91%  "
$ws.Range("E16").Value = "  -0.04%  "
$ws.Range("E17").Value = "  -0.29%  "
$ws.Range("D18").Value = "25.761.74"
$ws.Range("E18").Value = "  -0.49%  "
$ws.Range("D19").Value = "'1.00"
$ws.Range("D20").Value = "'4.43"
$ws.Range("E20").Value = "  +1.18%  "
$ws.Range("D21").Value = "'193.79"
$ws.Range("E21").Value = "  -0.17%  "
$ws.Range("E22").Value = "  -0.28%  "
$ws.Range("E23").Value = "  +1.65%  "
$ws.Range("D24").Value = "'1.00"
$ws.Range("E24").Value = "  +0.17%  "
$ws.Range("E25").Value = "  -0.66%  "
$ws.Range("D26").Value = "'140.16"
$ws.Range("E26").Value = "  -0.17%  "
$ws.Range("E27").Value = "  -3.94%  "
$ws.Range("E28").Value = "  -0.35%  "
$ws.Range("D29").Value = "'15.51"
$ws.Range("E29").Value = "  -0.36%  "
$ws.Range("E30").Value = "  -0.59%  "
$ws.Range("E31").Value = "  -2.16%  "
$ws.Range("D32").Value = "'3.34"
$ws.Range("E32").Value = "  +0.63%  "
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("E34").Value = "  +0.11%  "
$ws.Range("E35").Value = "  +0.37%  "
$ws.Range("D36").Value = "'0.896"
$ws.Range("E36").Value = "  -1.41%  "
$ws.Range("D37").Value = "'2.56"
$ws.Range("D38").Value = "'0.548"
$ws.Range("E38").Value = "  -1.67%  "
$ws.Range("D39").Value = "1.106.08"
$ws.Range("E39").Value = "  -2.20%  "
$ws.Range("D40").Value = "'0.0156"
$ws.Range("E40").Value = "  -0.75%  "
$ws.Range("E41").Value = "  +0.09%  "
$ws.Range("E42").Value = "  +0.77%  "
$ws.Range("D43").Value = "'99.91"
$ws.Range("E43").Value = "  +1.40%  "
$ws.Range("D44").Value = "'0.802"
$ws.Range("E44").Value = "  -0.27%  "
$ws.Range("D45").Value = "0.0₆0108"
$ws.Range("E45").Value = "  -4.75%  "
$ws.Range("D46").Value = "'55.04"
$ws.Range("E46").Value = "  -1.49%  "
$ws.Range("E47").Value = "  -1.97%  "
$ws.Range("D48").Value = "'7.71"
$ws.Range("E48").Value = "  -1.07%  "
$ws.Range("E49").Value = "  -0.26%  "
$ws.Range("E50").Value = "  +3.61%  "
$ws.Range("E51").Value = "  +0.38%  "

Write-Host "Applied cryptos update"
